# Insert two new data rows at the top of the "Vega Modelo de Temuco - Coliflor"
# block (rows 431-536), pushing the existing 431-536 rows down to 433-538,
# then populate the two newly-inserted rows (431-432) with their new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 431:536 down by two rows -> they become 433:538.
$ws.Rows("431:432").Insert()

# New row 431
$ws.Range("A431").Value = 10
$ws.Range("B431").Value = "Vega Modelo de Temuco"
$ws.Range("C431").Value = "La Araucanía"
$ws.Range("D431").Value = 44943
$ws.Range("E431").Value = 9
$ws.Range("F431").Value = 100112008
$ws.Range("G431").Value = "Coliflor"
$ws.Range("H431").Value = "Sin especificar"
$ws.Range("I431").Value = "Primera"
$ws.Range("J431").Value = 750
$ws.Range("K431").Value = 1000
$ws.Range("L431").Value = 1000
$ws.Range("M431").Value = 1000
$ws.Range("N431").Value = "$/unidad"
$ws.Range("O431").Value = "Región Metropolitana"
$ws.Range("P431").Value = 1000
$ws.Range("Q431").Value = 1
$ws.Range("R431").Value = "Hortaliza"

# New row 432
$ws.Range("A432").Value = 10
$ws.Range("B432").Value = "Vega Modelo de Temuco"
$ws.Range("C432").Value = "La Araucanía"
$ws.Range("D432").Value = 44943
$ws.Range("E432").Value = 9
$ws.Range("F432").Value = 100112008
$ws.Range("G432").Value = "Coliflor"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 1250
$ws.Range("K432").Value = 1000
$ws.Range("L432").Value = 1000
$ws.Range("M432").Value = 1000
$ws.Range("N432").Value = "$/unidad"
$ws.Range("O432").Value = "Región del Maule"
$ws.Range("P432").Value = 1000
$ws.Range("Q432").Value = 1
$ws.Range("R432").Value = "Hortaliza"
